$wb = $excel.ActiveWorkbook

$wsExh = $wb.Worksheets.Item("展览")
$wsAll = $wb.Worksheets.Item("全部类型")

# Sheet "展览" updates
$wsExh.Range("F5").Value = 2824
$wsExh.Range("F10").Value = 81
$wsExh.Range("F11").Value = 98
$wsExh.Range("F12").Value = 2657
$wsExh.Range("F13").Value = 854

# Sheet "全部类型" updates
$wsAll.Range("F6").Value = 2824
$wsAll.Range("F12").Value = 81
$wsAll.Range("F13").Value = 98
$wsAll.Range("F14").Value = 2657
$wsAll.Range("F15").Value = 854

$wb.Save()
